$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leetcode")

# Row 35 - Median of Two Sorted Arrays (Deferred)
$ws.Range("A35").Value = "Leetcode"
$ws.Range("B35").Value = 4
$ws.Range("D35").Value = "Arrays, Binary Search"
$ws.Range("E35").Value = "Hard"
$ws.Range("F35").Value = "Neetcode 150"
$ws.Range("C35").Value = "Median of Two Sorted Arrays"
$ws.Range("G35").Value = "DEFERRED"
$ws.Range("I35").Value = "Sticking with Easy and Medium for now."

# Row 36 - Reverse Linked List
$ws.Range("A36").Value = "Leetcode"
$ws.Range("B36").Value = 206
$ws.Range("E36").Value = "Easy"
$ws.Range("F36").Value = "Neetcode 150"
$ws.Range("G36").Value = "STRUGGLED"
$ws.Range("H36").Value = "17/06/2025"
$ws.Range("C36").Value = "Reverse Linked List"
$ws.Range("I36").Value = "Linked lists…I'm scared!"
$ws.Range("D36").Value = "Linked Lists, Two Pointers, Recursion"

# Row 37 - Merge Two Sorted Lists
$ws.Range("A37").Value = "Leetcode"
$ws.Range("B37").Value = 21
$ws.Range("E37").Value = "Easy"
$ws.Range("G37").Value = "STRUGGLED"
$ws.Range("C37").Value = "Merge Two Sorted Lists"
$ws.Range("D37").Value = "Linked Lists"
$ws.Range("F37").Value = "Neetcod 150"
$ws.Range("H37").Value = "18/06/2025"
$ws.Range("I37").Value = "Making a bit more sense."

# Mirror the end-user's final selection/scroll state (frozen header row is
# left untouched so the existing ySplit="1" freeze survives the round-trip).
$ws.Activate()
$ws.Range("G38").Select()
